$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 7303
$ws.Range("J3").Value = 7682
$ws.Range("H4").Value = 1712
$ws.Range("J4").Value = 1670
$ws.Range("J5").Value = 603
$ws.Range("J6").Value = 10477
$ws.Range("H7").Value = 26023
$ws.Range("J7").Value = 27735

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 272
$ws.Range("J7").Value = 416

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J2").Value = 461
$ws.Range("J3").Value = 508
$ws.Range("J6").Value = 644
$ws.Range("J7").Value = 1748

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J3").Value = 205
$ws.Range("J7").Value = 554

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 286
$ws.Range("J3").Value = 417
$ws.Range("J7").Value = 1258

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J3").Value = 143
$ws.Range("J7").Value = 397

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 256
$ws.Range("J3").Value = 284
$ws.Range("J6").Value = 248
$ws.Range("J7").Value = 850

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 201
$ws.Range("J6").Value = 261
$ws.Range("J7").Value = 699

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 219
$ws.Range("J4").Value = 128
$ws.Range("J7").Value = 793
$ws.Range("J8").Value = 1748
$ws.Range("J12").Value = 55
$ws.Range("J18").Value = 224
$ws.Range("J19").Value = 796
$ws.Range("J20").Value = 595
$ws.Range("J23").Value = 255
$ws.Range("J24").Value = 92
$ws.Range("J27").Value = 166
$ws.Range("J29").Value = 1478
$ws.Range("J31").Value = 288
$ws.Range("J33").Value = 1258
$ws.Range("J36").Value = 373
$ws.Range("J37").Value = 850
$ws.Range("J40").Value = 62
$ws.Range("J42").Value = 1186
$ws.Range("J43").Value = 234
$ws.Range("J44").Value = 219
$ws.Range("J47").Value = 202
$ws.Range("J51").Value = 343
$ws.Range("J52").Value = 710
$ws.Range("J53").Value = 416
$ws.Range("J54").Value = 549
$ws.Range("J60").Value = 164
$ws.Range("H63").Value = 270
$ws.Range("J63").Value = 82
$ws.Range("J65").Value = 699
$ws.Range("J66").Value = 84
$ws.Range("J73").Value = 268
$ws.Range("J77").Value = 194
$ws.Range("J79").Value = 764
$ws.Range("J80").Value = 50
$ws.Range("J83").Value = 554
$ws.Range("J84").Value = 231
$ws.Range("J85").Value = 1135
$ws.Range("J88").Value = 294
$ws.Range("J89").Value = 343
$ws.Range("J90").Value = 292
$ws.Range("J91").Value = 318
$ws.Range("J93").Value = 120
$ws.Range("J94").Value = 307
$ws.Range("J95").Value = 397
$ws.Range("H101").Value = 26023
$ws.Range("J101").Value = 27735

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 288

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J2").Value = 68
$ws.Range("J6").Value = 76
$ws.Range("J7").Value = 231

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 254
$ws.Range("J7").Value = 549

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 523
$ws.Range("J6").Value = 374
$ws.Range("J7").Value = 1478

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 195
$ws.Range("J3").Value = 228
$ws.Range("J7").Value = 796

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J6").Value = 88
$ws.Range("J7").Value = 219

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J6").Value = 628
$ws.Range("J7").Value = 1186

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("J3").Value = 24
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("J3").Value = 84
$ws.Range("J7").Value = 255

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 83
$ws.Range("J6").Value = 82
$ws.Range("J7").Value = 318

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J2").Value = 216
$ws.Range("J3").Value = 255
$ws.Range("J6").Value = 230
$ws.Range("J7").Value = 764

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J4").Value = 46
$ws.Range("J5").Value = 16
$ws.Range("J6").Value = 173
$ws.Range("J7").Value = 595

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 224

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J2").Value = 122
$ws.Range("J7").Value = 373

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 120

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 250
$ws.Range("J5").Value = 21
$ws.Range("J7").Value = 793

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J2").Value = 60
$ws.Range("J6").Value = 162
$ws.Range("J7").Value = 307

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 202

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 84

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J3").Value = 66
$ws.Range("J7").Value = 268

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J2").Value = 66
$ws.Range("J7").Value = 219

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 155
$ws.Range("J7").Value = 294

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 103
$ws.Range("J7").Value = 343

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J4").Value = 20
$ws.Range("J7").Value = 166

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J2").Value = 104
$ws.Range("J7").Value = 292

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("J3").Value = 90
$ws.Range("J6").Value = 141
$ws.Range("J7").Value = 343

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J3").Value = 45
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J6").Value = 139
$ws.Range("J7").Value = 234

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 305
$ws.Range("J7").Value = 1135

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J3").Value = 63
$ws.Range("J6").Value = 35
$ws.Range("J7").Value = 194

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 50

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("J3").Value = 22
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 167
$ws.Range("J3").Value = 198
$ws.Range("J7").Value = 710

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J2").Value = 39
$ws.Range("J7").Value = 128

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("J6").Value = 33
$ws.Range("J7").Value = 55
